# Update the "取得日時" (acquisition timestamp) column on the ランサーズ sheet
# from 2025-11-21 12:35:38 to 2025-11-21 12:46:41 for all data rows (2-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-21 12:46:41"

for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
